$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2597864768683274
$ws1.Range("C2").Value = 0.06108597285067873
$ws1.Range("D2").Value = 0.9642857142857143
$ws1.Range("E2").Value = 0.1148936170212766
$ws1.Range("F2").Value = 0.2436823104693141
$ws1.Range("G2").Value = 0.6147110332749562
$ws1.Range("H2").Value = 0.8214620117710005
$ws1.Range("I2").Value = 27
$ws1.Range("J2").Value = 415
$ws1.Range("K2").Value = 119
$ws1.Range("L2").Value = 1

# --- Classification Report sheet ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$ws2.Range("B2").Value = 0.9916666666666667
$ws2.Range("C2").Value = 0.2228464419475655
$ws2.Range("D2").Value = 0.363914373088685

# Row 3 ("1")
$ws2.Range("B3").Value = 0.06108597285067873
$ws2.Range("C3").Value = 0.9642857142857143
$ws2.Range("D3").Value = 0.1148936170212766

# Row 4 ("accuracy")
$ws2.Range("B4").Value = 0.2597864768683274
$ws2.Range("C4").Value = 0.2597864768683274
$ws2.Range("D4").Value = 0.2597864768683274
$ws2.Range("E4").Value = 0.2597864768683274

# Row 5 ("macro avg")
$ws2.Range("B5").Value = 0.5263763197586727
$ws2.Range("C5").Value = 0.5935660781166399
$ws2.Range("D5").Value = 0.2394039950549808

# Row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.9453032157292154
$ws2.Range("C6").Value = 0.2597864768683274
$ws2.Range("D6").Value = 0.3515076450283871

# --- Confusion Matrix sheet ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$ws3.Range("B2").Value = 119
$ws3.Range("C2").Value = 415

# Row 3 ("Actual 1")
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 27
